$d = $word.ActiveDocument

# --- locate the "react-moment" paragraph; the new paragraph goes right after it ---
$anchorRng = $d.Content
$anchorRng.Find.Execute("react-moment", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorPara = $anchorRng.Paragraphs(1)
$insertAt = $anchorPara.Range.End

# --- locate an existing "npm install <pkg>" paragraph to use as a formatting template,
#     so the new paragraph keeps the same run layout (word / " install " / word) ---
$templateRng = $d.Content
$templateRng.Find.Execute("uuid", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$templatePara = $templateRng.Paragraphs(1)
$templateFullRng = $d.Range($templatePara.Range.Start, $templatePara.Range.End)

# Copy the whole template paragraph (incl. its paragraph mark) and paste it right
# after the "react-moment" paragraph, pushing the rest of the document down.
$templateFullRng.Copy()
$destRng = $d.Range($insertAt, $insertAt)
$destRng.Paste()

# --- the newly pasted paragraph now sits right after the "react-moment" paragraph ---
$newPara = $anchorPara.Next()
$newParaEnd = $newPara.Range.End

# Replace the trailing "uuid" word with "swiper" without disturbing the other runs:
# delete just the word, then insert the new word in its place.
$wordEnd = $newParaEnd - 1
$wordStart = $wordEnd - 4
$wordRng = $d.Range($wordStart, $wordEnd)
$wordRng.Delete()
$insPoint = $d.Range($wordStart, $wordStart)
$insPoint.InsertAfter("swiper")
